# Katalog guncellendi - Cum 14.11.2025 11:32:48,27
# Appends extra "satin alma" sentences to the three shared product
# descriptions (aciklama, column E) and restores the last on-screen
# selection (E28) that the author left the workbook in.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$suffixTekli  = "Ürünümüz tekli olarak satın alınabilir."
$suffixSerili = "Ürünümüz serili olarak satılmaktadır."

# Rows 26-30: "Eşref Gömlek *" shirts (Polyester materyali ... description)
for ($r = 26; $r -le 30; $r++) {
    $cell = $ws.Cells.Item($r, 5)
    $cell.Value = $cell.Value2 + $suffixSerili
}

# Rows 23-25: "Kot Gömlek Ceket *" shirt/jackets (%100 pamuk ... description)
for ($r = 23; $r -le 25; $r++) {
    $cell = $ws.Cells.Item($r, 5)
    $cell.Value = $cell.Value2 + $suffixTekli
}

# Rows 11-22: "Kot Gömlek *" shirts (%85 pamuk ... description)
for ($r = 11; $r -le 22; $r++) {
    $cell = $ws.Cells.Item($r, 5)
    $cell.Value = $cell.Value2 + $suffixTekli
}

# Restore the author's last selection/scroll position in the sheet.
$ws.Range("E28").Select()
